# Commit automatique via PowerShell
# - Replace "République Démocratique du Congo" with "RDC" in the PAYS column
# - Narrow columns B (NOM) and E (PAYS) now that the country name is shorter

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update country name for every row that still references the long form.
$ws.Range("E14").Value = "RDC"
$ws.Range("E22").Value = "RDC"
$ws.Range("E23").Value = "RDC"
$ws.Range("E24").Value = "RDC"
$ws.Range("E25").Value = "RDC"

# Shrink column widths to fit the now much shorter "RDC" text.
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.666666666666666
